# Patient log update: add newly logged UIC cases (rows 779-810).
# Mirrors the author's commit "updated the name and journal club".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting for the MRN (B, text) and DATE (C, date) columns
# --- from the last pre-existing data row (778) onto the new rows, so the
# --- new cells reuse the same style indexes instead of creating new ones.
$ws.Range("B778").Copy() | Out-Null
$ws.Range("B779:B810").PasteSpecial(-4122) | Out-Null
$ws.Range("C778").Copy() | Out-Null
$ws.Range("C779:C810").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Populate the new rows. They are written in the same order the data was
# --- originally entered (779-804, then 806-810, then 805 last) so that the
# --- shared-string table is appended to in the same order as the source file.
# Row 779
$ws.Range("A779").Value = 'UIC'
$ws.Range("B779").Value = '080936335'
$ws.Range("C779").Value = 44223
$ws.Range("G779").Value = 'TEE'
$ws.Range("H779").Value = 'bedside'

# Row 780
$ws.Range("A780").Value = 'UIC'
$ws.Range("B780").Value = '080117366'
$ws.Range("C780").Value = 44223
$ws.Range("G780").Value = 'TTE'

# Row 781
$ws.Range("A781").Value = 'UIC'
$ws.Range("B781").Value = '080021471'
$ws.Range("C781").Value = 44223
$ws.Range("G781").Value = 'TTE'

# Row 782
$ws.Range("A782").Value = 'UIC'
$ws.Range("B782").Value = '081155573'
$ws.Range("C782").Value = 44223
$ws.Range("G782").Value = 'stress'
$ws.Range("H782").Value = 'treadmill'
$ws.Range("I782").Value = 'echo'

# Row 783
$ws.Range("A783").Value = 'UIC'
$ws.Range("B783").Value = '081450956'
$ws.Range("C783").Value = 44223
$ws.Range("G783").Value = 'TTE'

# Row 784
$ws.Range("A784").Value = 'UIC'
$ws.Range("B784").Value = '080936335'
$ws.Range("C784").Value = 44223
$ws.Range("G784").Value = 'TEE'

# Row 785
$ws.Range("A785").Value = 'UIC'
$ws.Range("B785").Value = '081482602'
$ws.Range("C785").Value = 44223
$ws.Range("G785").Value = 'TTE'
$ws.Range("H785").Value = 'bedside'

# Row 786
$ws.Range("A786").Value = 'UIC'
$ws.Range("B786").Value = '200253384'
$ws.Range("C786").Value = 44223
$ws.Range("G786").Value = 'TTE'

# Row 787
$ws.Range("A787").Value = 'UIC'
$ws.Range("B787").Value = '081520194'
$ws.Range("C787").Value = 44222
$ws.Range("G787").Value = 'TTE'

# Row 788
$ws.Range("A788").Value = 'UIC'
$ws.Range("B788").Value = '080655007'
$ws.Range("C788").Value = 44222
$ws.Range("G788").Value = 'TTE'

# Row 789
$ws.Range("A789").Value = 'UIC'
$ws.Range("B789").Value = '081482161'
$ws.Range("C789").Value = 44222
$ws.Range("G789").Value = 'TTE'
$ws.Range("H789").Value = 'bedside'

# Row 790
$ws.Range("A790").Value = 'UIC'
$ws.Range("B790").Value = '051641983'
$ws.Range("C790").Value = 44222
$ws.Range("G790").Value = 'TTE'

# Row 791
$ws.Range("A791").Value = 'UIC'
$ws.Range("B791").Value = '080678498'
$ws.Range("C791").Value = 44222
$ws.Range("G791").Value = 'TTE'

# Row 792
$ws.Range("A792").Value = 'UIC'
$ws.Range("B792").Value = '200256419'
$ws.Range("C792").Value = 44222
$ws.Range("G792").Value = 'TTE'

# Row 793
$ws.Range("A793").Value = 'UIC'
$ws.Range("B793").Value = '081323814'
$ws.Range("C793").Value = 44222
$ws.Range("G793").Value = 'TTE'

# Row 794
$ws.Range("A794").Value = 'UIC'
$ws.Range("B794").Value = '081482602'
$ws.Range("C794").Value = 44222
$ws.Range("D794").Value = 'RHD'
$ws.Range("E794").Value = 'MS'
$ws.Range("F794").Value = 'AS'
$ws.Range("G794").Value = 'TEE'
$ws.Range("H794").Value = 'bedside'
$ws.Range("I794").Value = 'structural'

# Row 795
$ws.Range("A795").Value = 'UIC'
$ws.Range("B795").Value = '200195082'
$ws.Range("C795").Value = 44221
$ws.Range("G795").Value = 'TTE'

# Row 796
$ws.Range("A796").Value = 'UIC'
$ws.Range("B796").Value = '081548675'
$ws.Range("C796").Value = 44221
$ws.Range("G796").Value = 'TTE'

# Row 797
$ws.Range("A797").Value = 'UIC'
$ws.Range("B797").Value = '080310420'
$ws.Range("C797").Value = 44221
$ws.Range("G797").Value = 'TTE'

# Row 798
$ws.Range("A798").Value = 'UIC'
$ws.Range("B798").Value = '200106873'
$ws.Range("C798").Value = 44221
$ws.Range("G798").Value = 'TTE'

# Row 799
$ws.Range("A799").Value = 'UIC'
$ws.Range("B799").Value = '080878880'
$ws.Range("C799").Value = 44221
$ws.Range("G799").Value = 'TTE'

# Row 800
$ws.Range("A800").Value = 'UIC'
$ws.Range("B800").Value = '072952948'
$ws.Range("C800").Value = 44221
$ws.Range("G800").Value = 'TTE'

# Row 801
$ws.Range("A801").Value = 'UIC'
$ws.Range("B801").Value = '200255391'
$ws.Range("C801").Value = 44221
$ws.Range("G801").Value = 'TTE'

# Row 802
$ws.Range("A802").Value = 'UIC'
$ws.Range("B802").Value = '006640692'
$ws.Range("C802").Value = 44221
$ws.Range("G802").Value = 'TTE'

# Row 803
$ws.Range("A803").Value = 'UIC'
$ws.Range("B803").Value = '073052763'
$ws.Range("C803").Value = 44221
$ws.Range("G803").Value = 'TTE'

# Row 804
$ws.Range("A804").Value = 'UIC'
$ws.Range("B804").Value = '050988922'
$ws.Range("C804").Value = 44221
$ws.Range("G804").Value = 'TTE'

# Row 806
$ws.Range("A806").Value = 'UIC'
$ws.Range("B806").Value = '074413352'
$ws.Range("C806").Value = 44221
$ws.Range("G806").Value = 'TTE'

# Row 807
$ws.Range("A807").Value = 'UIC'
$ws.Range("B807").Value = '081551473'
$ws.Range("C807").Value = 44221
$ws.Range("D807").Value = 'PFO'
$ws.Range("E807").Value = 'CVA'
$ws.Range("G807").Value = 'TEE'
$ws.Range("H807").Value = 'bedside'

# Row 808
$ws.Range("A808").Value = 'UIC'
$ws.Range("B808").Value = '081551473'
$ws.Range("C808").Value = 44221
$ws.Range("D808").Value = 'PFO'
$ws.Range("E808").Value = 'CVA'
$ws.Range("G808").Value = 'TEE'

# Row 809
$ws.Range("A809").Value = 'UIC'
$ws.Range("B809").Value = '080606193'
$ws.Range("C809").Value = 44221
$ws.Range("G809").Value = 'TTE'

# Row 810
$ws.Range("A810").Value = 'UIC'
$ws.Range("B810").Value = '200255583'
$ws.Range("C810").Value = 44221
$ws.Range("G810").Value = 'TTE'

# Row 805
$ws.Range("A805").Value = 'UIC'
$ws.Range("B805").Value = '031295212'
$ws.Range("C805").Value = 44221
$ws.Range("G805").Value = 'TTE'

# --- Restore the active-cell selection to match the edited view state. ---
$ws.Activate() | Out-Null
$ws.Range("D807").Select() | Out-Null

Write-Output "Added rows 779-810; sharedStrings/sheet now in sync with log update."
